# Regenerate the handback status report:
#  - the "ea42df71-...-ae60a86de51d" entry (row 3) has been excluded/handled, so
#    its row is dropped from all three sheets.
#  - the remaining "79efc3a7-..." entry's handoff/handback timestamps are refreshed
#    to the new report run.

$wb = $excel.ActiveWorkbook

function Remove-HyperlinksByAddress($ws, $addresses) {
    $toDelete = @()
    foreach ($h in $ws.Hyperlinks) {
        if ($addresses -contains $h.Range.Address()) {
            $toDelete += $h
        }
    }
    for ($i = $toDelete.Count - 1; $i -ge 0; $i--) {
        $toDelete[$i].Delete()
    }
}

# --- Overview sheet -------------------------------------------------------
$ovw = $wb.Worksheets.Item("Overview")
Remove-HyperlinksByAddress $ovw @('$A$3')
$ovw.Rows("3:3").Delete()

# --- zh-cn sheet ------------------------------------------------------------
$zh = $wb.Worksheets.Item("zh-cn")
Remove-HyperlinksByAddress $zh @('$A$3','$B$3','$D$3','$F$3','$G$3')
$zh.Rows("3:3").Delete()
$zh.Range("E2").Value = "2016-03-19 02:38:19"
$zh.Range("H2").Value = "2016-03-19 02:38:37"

# --- de-de sheet ------------------------------------------------------------
$de = $wb.Worksheets.Item("de-de")
Remove-HyperlinksByAddress $de @('$A$3','$B$3','$D$3','$F$3','$G$3')
$de.Rows("3:3").Delete()
$de.Range("E2").Value = "2016-03-19 02:38:21"
$de.Range("H2").Value = "2016-03-19 02:38:42"
